# Fixed update to excel issue
#
# - Rename the "Requested quantity" header on the "Weekly Quantity" sheet to
#   "Weekly_PO_Qty".
# - Rename the "Requested quantity" header on the "Monthly Trend" sheet to
#   "Monthly_PO_Qty".
# - Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)
#   with the forecasted weekly PO quantities.

$wb = $excel.ActiveWorkbook

$weekly  = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# --- header renames -------------------------------------------------------
$weekly.Range("B1").Value  = "Weekly_PO_Qty"
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- new "PO Forecast" sheet ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fc = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$fc.Name = "PO Forecast"

# Match the outline / page-margin conventions used by the other sheets.
$fc.Outline.SummaryRow    = 1
$fc.Outline.SummaryColumn = 1
$fc.PageSetup.LeftMargin   = 0.75 * 72
$fc.PageSetup.RightMargin  = 0.75 * 72
$fc.PageSetup.TopMargin    = 1    * 72
$fc.PageSetup.BottomMargin = 1    * 72
$fc.PageSetup.HeaderMargin = 0.5  * 72
$fc.PageSetup.FooterMargin = 0.5  * 72

$fc.Range("A1").Value = "ds"
$fc.Range("B1").Value = "PO_Forecast"
$fc.Range("C1").Value = "yhat_lower"
$fc.Range("D1").Value = "yhat_upper"

# Copy the header formatting (bold, centered, bordered) from the existing
# sheets so the new header reuses the same cell style instead of minting a
# new one.
$weekly.Range("A1:B1").Copy()
$fc.Range("A1:D1").PasteSpecial(-4122)

$data = @(
    @(45312.99999999999,71,-274.489409555782,396.067746716375),
    @(45319.99999999999,77,-281.2128042529649,413.267553057409),
    @(45403.99999999999,148,-186.6172412035,481.8538057115864),
    @(45410.99999999999,154,-182.6953230288297,466.4545780367071),
    @(45424.99999999999,165,-170.8890051239006,514.3710230808538),
    @(45445.99999999999,183,-174.8843750267436,517.3045967221565),
    @(45452.99999999999,189,-157.3855896010891,540.3502709302746),
    @(45466.99999999999,201,-139.0939473023811,520.0706594768222),
    @(45473.99999999999,207,-162.557576454658,550.1891118956687),
    @(45480.99999999999,213,-135.3668555235642,540.2855241969396),
    @(45487.99999999999,219,-112.9568064971525,554.2418591196105),
    @(45522.99999999999,248,-57.07718015950802,592.8884285934466),
    @(45529.99999999999,254,-69.09606649411907,606.8410252654982),
    @(45536.99999999999,260,-49.06836385626646,615.2922163472905),
    @(45564.99999999999,284,-52.73867601724429,602.0149203809727),
    @(45571.99999999999,289,-63.50987958238752,614.892595547775),
    @(45585.99999999999,301,-26.85159084392301,616.1124434099165),
    @(45599.99999999999,313,-20.30473190515999,667.403837429843),
    @(45613.99999999999,325,8.949547358030813,663.1671996149353),
    @(45620.99999999999,331,-9.853081831836755,675.5490765985616),
    @(45627.99999999999,337,-6.800000077947697,665.6988776245722),
    @(45634.99999999999,343,21.75552979021066,699.9887418448557),
    @(45641.99999999999,349,1.409164462247696,659.6769184322196),
    @(45648.99999999999,355,27.6324397506751,683.045356730385),
    @(45655.99999999999,360,-6.461771780183521,695.7947591981757),
    @(45662.99999999999,366,27.07420583872602,708.4068854598256),
    @(45669.99999999999,372,44.16751989812407,701.5125556508866)
)

$r = 2
foreach ($row in $data) {
    $fc.Cells.Item($r, 1).Value = $row[0]
    $fc.Cells.Item($r, 2).Value = $row[1]
    $fc.Cells.Item($r, 3).Value = $row[2]
    $fc.Cells.Item($r, 4).Value = $row[3]
    $r++
}
$lastDataRow = $r - 1

# Copy the date/time number-format used for column A on the other sheets.
$weekly.Range("A2").Copy()
$fc.Range("A2:A$lastDataRow").PasteSpecial(-4122)
